# Daily refresh of the cryptos price sheet (GitHub Actions job), updating
# coin prices, swapping in newly-ranked coins (WazirX moves up into row 9,
# pushing BitrueCoin/BitMartToken/BitForexToken/One/TigerCash/LEO/GateToken/
# BTSEToken/BitpandaEcosystemToken down a row; KickToken/CEJI swap places)
# and refreshing the "Volume(1h)" labels to match.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Non-numeric text cells (Coin names, URLs, labels) - plain assignment is safe
$ws.Range("B9").Value = 'WazirX'
$ws.Range("C9").Value = 'https://coinranking.com/coin/6QK-8hUZ+wazirx-wrx'
$ws.Range("E9").Value = '8WazirXWRX'
$ws.Range("B10").Value = 'BitrueCoin'
$ws.Range("C10").Value = 'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr'
$ws.Range("E10").Value = '9BitrueCoinBTR'
$ws.Range("B11").Value = 'BitMartToken'
$ws.Range("C11").Value = 'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx'
$ws.Range("E11").Value = '10BitMartTokenBMX'
$ws.Range("B12").Value = 'BitForexToken'
$ws.Range("C12").Value = 'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf'
$ws.Range("E12").Value = '11BitForexTokenBF'
$ws.Range("B13").Value = 'One'
$ws.Range("C13").Value = 'https://coinranking.com/coin/6Lga5NiXX3rT+one-one'
$ws.Range("E13").Value = '12OneONEBestin24h'
$ws.Range("B14").Value = 'TigerCash'
$ws.Range("C14").Value = 'https://coinranking.com/coin/6hIn06L2+tigercash-tch'
$ws.Range("E14").Value = '13TigerCashTCH'
$ws.Range("B15").Value = 'LEO'
$ws.Range("C15").Value = 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'
$ws.Range("E15").Value = '14LEOLEO'
$ws.Range("B16").Value = 'GateToken'
$ws.Range("C16").Value = 'https://coinranking.com/coin/t7m8DZVyMsAu+gatetoken-gt'
$ws.Range("E16").Value = '15GateTokenGT'
$ws.Range("B17").Value = 'BTSEToken'
$ws.Range("C17").Value = 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'
$ws.Range("E17").Value = '16BTSETokenBTSE'
$ws.Range("B18").Value = 'BitpandaEcosystemToken'
$ws.Range("C18").Value = 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'
$ws.Range("E18").Value = '17BitpandaEcosystemTokenBEST'
$ws.Range("E41").Value = '40BKEXTokenBKK'
$ws.Range("B42").Value = 'KickToken'
$ws.Range("C42").Value = 'https://coinranking.com/coin/F_Yv9Cu7pPL3Y+kicktoken-kick'
$ws.Range("E42").Value = '41KickTokenKICK'
$ws.Range("B43").Value = 'CEJI'
$ws.Range("C43").Value = 'https://coinranking.com/coin/SbKjCVJCh+ceji-ceji'
$ws.Range("E43").Value = '42CEJICEJI'

# Numeric-looking text cells (Price column) - force text to preserve exact
# formatting (leading/trailing zeros) without leaving a quote-prefix style behind
$numericTextCells = @{
    'D2' = '246.66'
    'D3' = '26.54'
    'D4' = '5.085'
    'D5' = '0.05614'
    'D6' = '6.494'
    'D7' = '0.8138'
    'D8' = '0.8455'
    'D9' = '0.1343'
    'D10' = '0.02852'
    'D11' = '0.09392'
    'D12' = '0.001533'
    'D13' = '0.009904'
    'D14' = '0.006159'
    'D15' = '3.588'
    'D16' = '3.012'
    'D17' = '2.118'
    'D18' = '0.3157'
    'D19' = '0.06965'
    'D20' = '0.03200'
    'D21' = '0.1321'
    'D22' = '3.740'
    'D23' = '0.04653'
    'D24' = '0.1351'
    'D25' = '0.001249'
    'D26' = '0.004590'
    'D27' = '0.00009604'
    'D28' = '0.0001939'
    'D40' = '0.03670'
    'D41' = '0.1057'
    'D42' = '0.006183'
    'D43' = '0.002604'
    'D44' = '0.008894'
    'D45' = '0.00005299'
    'D47' = '0.1201'
    'D48' = '0.002519'
    'D49' = '0.00002101'
    'D50' = '0.0002001'
}
foreach ($addr in $numericTextCells.Keys) {
    $cell = $ws.Range($addr)
    $cell.NumberFormat = "@"
    $cell.Value = $numericTextCells[$addr]
    $cell.Style = "Normal"
}
